$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell (D1) + border under the whole header row ---
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("D1").Value = "Observaciones"
$ws.Range("A1:D1").Borders.Item(9).LineStyle = 1
$ws.Range("A1:D1").Borders.Item(9).Weight = 2

# --- Extend formatting for the new rows 12:16 (copy format from row 11) ---
$ws.Range("A11:C11").Copy()
$ws.Range("A12:C16").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# --- Values for existing rows that gain new annotation columns (D/E) ---
$ws.Range("E2").Value = "Profiling"
$ws.Range("E3").Value = "Profiling"
$ws.Range("E4").Value = "Profiling"
$ws.Range("D5").Value = "Limitación de 1000 peticiones al día y 50 al segundo"
$ws.Range("E5").Value = "Profiling?"
$ws.Range("D6").Value = "Limitación de 1000 peticiones al día y 50 al segundo"
$ws.Range("E6").Value = "Profiling?"

# --- New rows 12:16 ---
$ws.Range("A12").Value = "RegistrarMedicamento"
$ws.Range("B12").Value = 4200
$ws.Range("C12").Value = 5000

$ws.Range("A13").Value = "ActualizarMedicamento"
$ws.Range("B13").Value = 7000
$ws.Range("C13").Value = 9000

$ws.Range("A14").Value = "RegistrarMedialTest"
$ws.Range("B14").Value = 500
$ws.Range("C14").Value = 3000
$ws.Range("E14").Value = "Profiling"

$ws.Range("A15").Value = "RegistrarPetType"
$ws.Range("B15").Value = 5000
$ws.Range("C15").Value = 7000

$ws.Range("A16").Value = "ActualizarPetType"
$ws.Range("B16").Value = 7000
$ws.Range("C16").Value = 8500

# --- Column width for the new "Observaciones" column ---
$ws.Columns.Item(4).ColumnWidth = 60.43

# --- Page setup (A4, portrait) ---
$ws.PageSetup.PaperSize = [Microsoft.Office.Interop.Excel.XlPaperSize]::xlPaperA4
$ws.PageSetup.Orientation = [Microsoft.Office.Interop.Excel.XlPageOrientation]::xlPortrait

# --- Selection moves to C17 after the new data ---
$ws.Range("C17").Select() | Out-Null
